# Update "想去人数" (F column) values across sheets to match
# regenerated site data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 201
$ws1.Range("F5").Value  = 981
$ws1.Range("F6").Value  = 5428
$ws1.Range("F7").Value  = 476
$ws1.Range("F8").Value  = 669
$ws1.Range("F13").Value = 581
$ws1.Range("F17").Value = 1811
$ws1.Range("F19").Value = 880
$ws1.Range("F21").Value = 192
$ws1.Range("F22").Value = 327
$ws1.Range("F23").Value = 532
$ws1.Range("F24").Value = 143
$ws1.Range("F28").Value = 2796
$ws1.Range("F32").Value = 117
$ws1.Range("F34").Value = 356
$ws1.Range("F39").Value = 284
$ws1.Range("F40").Value = 684

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 180

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 201
$ws4.Range("F5").Value  = 981
$ws4.Range("F7").Value  = 5428
$ws4.Range("F8").Value  = 476
$ws4.Range("F9").Value  = 669
$ws4.Range("F11").Value = 180
$ws4.Range("F18").Value = 581
$ws4.Range("F23").Value = 1811
$ws4.Range("F25").Value = 880
$ws4.Range("F26").Value = 192
$ws4.Range("F27").Value = 327
$ws4.Range("F29").Value = 532
$ws4.Range("F30").Value = 143
$ws4.Range("F32").Value = 2796
$ws4.Range("F36").Value = 117
$ws4.Range("F38").Value = 356
$ws4.Range("F42").Value = 284
$ws4.Range("F43").Value = 684

$wb.Save()
